$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 137, shifting existing rows 137:140 down to 138:141
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row 137 with the new weekly record
$ws.Cells.Item(137, 1).Value = 9
$ws.Cells.Item(137, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(137, 3).Value = "Metropolitana"
$ws.Cells.Item(137, 4).Value = 44448
$ws.Cells.Item(137, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(137, 5).Value = 13
$ws.Cells.Item(137, 6).Value = 300000001
$ws.Cells.Item(137, 7).Value = "Rabanito"
$ws.Cells.Item(137, 8).Value = "Sin especificar"
$ws.Cells.Item(137, 9).Value = "Primera"
$ws.Cells.Item(137, 10).Value = 7900
$ws.Cells.Item(137, 11).Value = 3500
$ws.Cells.Item(137, 12).Value = 4000
$ws.Cells.Item(137, 13).Value = 3747
$ws.Cells.Item(137, 14).Value = '$/cien unidades (volumen en unidades)'
$ws.Cells.Item(137, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(137, 16).Value = 37
$ws.Cells.Item(137, 17).Value = 100
$ws.Cells.Item(137, 18).Value = "Hortaliza"
